{"js": "// Highlight four agenda bullet lines (and their paragraph marks) in yellow,\n// matching the \"added code for gateway\" commit:\n//   - \"What is an API Gateway\"\n//   - \"Overview of Ocelot\"\n//   - \"Overview of Asynchronous Communication\"\n//   - \"Using RabbitMQ for Asynchronous Communication\"\nconst targets = [\n  \"What is an API Gateway\",\n  \"Overview of Ocelot\",\n  \"Overview of Asynchronous Communication\",\n  \"Using RabbitMQ for Asynchronous Communication\",\n];\n\nconst body = context.document.body;\n\nfor (const target of targets) {\n  const results = body.search(target, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    const found = results.items[i];\n\n    // Highlight the matched run text itself.\n    found.font.highlightColor = \"Yellow\";\n\n    // Also highlight the paragraph mark (the paragraph's own rPr), so the\n    // whole list-item paragraph (incl. its end-of-paragraph mark) is yellow,\n    // mirroring what Word does when you select & highlight a full line.\n    const para = found.paragraphs.getFirst();\n    para.font.highlightColor = \"Yellow\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Highlight four agenda bullet lines (and their paragraph marks) in yellow,\n# matching the \"added code for gateway\" commit:\n#   - \"What is an API Gateway\"\n#   - \"Overview of Ocelot\"\n#   - \"Overview of Asynchronous Communication\"\n#   - \"Using RabbitMQ for Asynchronous Communication\"\n\n$d = $word.ActiveDocument\n\n$targets = @(\n    \"What is an API Gateway\",\n    \"Overview of Ocelot\",\n    \"Overview of Asynchronous Communication\",\n    \"Using RabbitMQ for Asynchronous Communication\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    # Paragraph.Range.Text includes the trailing paragraph-mark character\n    # (Cr), so trim it before comparing against our plain-text targets.\n    $paraText = $p.Range.Text.TrimEnd([char]13)\n    foreach ($target in $targets) {\n        if ($paraText -eq $target) {\n            # Setting HighlightColorIndex via the Range's Font applies the\n            # highlight to both the run text and the paragraph mark (pPr/rPr),\n            # matching what Word does when the whole line is selected and\n            # highlighted.\n            $p.Range.Font.HighlightColorIndex = 7   # wdYellow\n        }\n    }\n}\n"}
